$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row9 = New-Object 'object[,]' 1,49
$row9[0,0] = "llama3:8b-instruct-q5_K_M"
$row9[0,1] = "llama3:70b"
$row9[0,2] = 100
$row9[0,3] = 400
$row9[0,4] = 2612.39
$row9[0,5] = 1290.58
$row9[0,6] = 0.2875
$row9[0,7] = "logs\llama3_8b_instruct_q5_K_M_llama3_70b_100_400_test_match.txt"
$row9[0,8] = 1290.58
$row9[0,9] = 0.94375
$row9[0,10] = "logs\llama3_8b_instruct_q5_K_M_llama3_70b_100_400_test_correct.txt"
$row9[0,11] = 1290.58
$row9[0,12] = 0.95625
$row9[0,13] = "logs\llama3_8b_instruct_q5_K_M_llama3_70b_100_400_test_executable.txt"
$row9[0,14] = 0
$row9[0,28] = "text"
$row9[0,29] = 140
$row9[0,30] = 0.1
$row9[0,31] = 150
$row9[0,32] = 0.9
$row9[0,33] = 5
$row9[0,34] = 1
$row9[0,35] = 1.2
$row9[0,36] = 1
$row9[0,37] = 1024
$row9[0,38] = 0
$row9[0,39] = 1321.81
$row9[0,40] = 0.26875
$row9[0,41] = "logs\llama3_8b_instruct_q5_K_M_llama3_70b_100_400_test_fewshot_match.txt"
$row9[0,42] = 1321.81
$row9[0,43] = 0.86875
$row9[0,44] = "logs\llama3_8b_instruct_q5_K_M_llama3_70b_100_400_test_fewshot_correct.txt"
$row9[0,45] = 1321.81
$row9[0,46] = 0.89375
$row9[0,47] = "logs\llama3_8b_instruct_q5_K_M_llama3_70b_100_400_test_fewshot_executable.txt"
$row9[0,48] = 0
$ws.Range("A9:AW9").Value = $row9

$row10 = New-Object 'object[,]' 1,49
$row10[0,0] = "command-r"
$row10[0,1] = "llama3:70b"
$row10[0,2] = 100
$row10[0,3] = 400
$row10[0,4] = 2709.84
$row10[0,5] = 1335.52
$row10[0,6] = 0.28125
$row10[0,7] = "logs\command_r_llama3_70b_100_400_test_match.txt"
$row10[0,8] = 1335.52
$row10[0,9] = 0.90625
$row10[0,10] = "logs\command_r_llama3_70b_100_400_test_correct.txt"
$row10[0,11] = 1335.52
$row10[0,12] = 0.94375
$row10[0,13] = "logs\command_r_llama3_70b_100_400_test_executable.txt"
$row10[0,14] = 0
$row10[0,28] = "text"
$row10[0,29] = 140
$row10[0,30] = 0.1
$row10[0,31] = 150
$row10[0,32] = 0.9
$row10[0,33] = 5
$row10[0,34] = 1
$row10[0,35] = 1.2
$row10[0,36] = 1
$row10[0,37] = 1024
$row10[0,38] = 0
$row10[0,39] = 1374.32
$row10[0,40] = 0.3125
$row10[0,41] = "logs\command_r_llama3_70b_100_400_test_fewshot_match.txt"
$row10[0,42] = 1374.32
$row10[0,43] = 0.8875
$row10[0,44] = "logs\command_r_llama3_70b_100_400_test_fewshot_correct.txt"
$row10[0,45] = 1374.32
$row10[0,46] = 0.925
$row10[0,47] = "logs\command_r_llama3_70b_100_400_test_fewshot_executable.txt"
$row10[0,48] = 0
$ws.Range("A10:AW10").Value = $row10

$row11 = New-Object 'object[,]' 1,49
$row11[0,0] = "aya:35b"
$row11[0,1] = "llama3:70b"
$row11[0,2] = 100
$row11[0,3] = 400
$row11[0,4] = 2695.8
$row11[0,5] = 1345.42
$row11[0,6] = 0.28125
$row11[0,7] = "logs\aya_35b_llama3_70b_100_400_test_match.txt"
$row11[0,8] = 1345.42
$row11[0,9] = 0.875
$row11[0,10] = "logs\aya_35b_llama3_70b_100_400_test_correct.txt"
$row11[0,11] = 1345.42
$row11[0,12] = 0.93125
$row11[0,13] = "logs\aya_35b_llama3_70b_100_400_test_executable.txt"
$row11[0,14] = 0
$row11[0,28] = "text"
$row11[0,29] = 140
$row11[0,30] = 0.1
$row11[0,31] = 150
$row11[0,32] = 0.9
$row11[0,33] = 5
$row11[0,34] = 1
$row11[0,35] = 1.2
$row11[0,36] = 1
$row11[0,37] = 1024
$row11[0,38] = 0
$row11[0,39] = 1350.38
$row11[0,40] = 0.2875
$row11[0,41] = "logs\aya_35b_llama3_70b_100_400_test_fewshot_match.txt"
$row11[0,42] = 1350.38
$row11[0,43] = 0.85
$row11[0,44] = "logs\aya_35b_llama3_70b_100_400_test_fewshot_correct.txt"
$row11[0,45] = 1350.38
$row11[0,46] = 0.9375
$row11[0,47] = "logs\aya_35b_llama3_70b_100_400_test_fewshot_executable.txt"
$row11[0,48] = 0
$ws.Range("A11:AW11").Value = $row11

$row12 = New-Object 'object[,]' 1,49
$row12[0,0] = "qwen2:7b-instruct-q5_K_M"
$row12[0,1] = "llama3:70b"
$row12[0,2] = 100
$row12[0,3] = 400
$row12[0,4] = 2632.75
$row12[0,5] = 1309.74
$row12[0,6] = 0.275
$row12[0,7] = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_100_400_test_match.txt"
$row12[0,8] = 1309.74
$row12[0,9] = 0.9
$row12[0,10] = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_100_400_test_correct.txt"
$row12[0,11] = 1309.74
$row12[0,12] = 0.93125
$row12[0,13] = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_100_400_test_executable.txt"
$row12[0,14] = 0
$row12[0,28] = "text"
$row12[0,29] = 140
$row12[0,30] = 0.1
$row12[0,31] = 150
$row12[0,32] = 0.9
$row12[0,33] = 5
$row12[0,34] = 1
$row12[0,35] = 1.2
$row12[0,36] = 1
$row12[0,37] = 1024
$row12[0,38] = 0
$row12[0,39] = 1323.01
$row12[0,40] = 0.3
$row12[0,41] = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_100_400_test_fewshot_match.txt"
$row12[0,42] = 1323.01
$row12[0,43] = 0.86875
$row12[0,44] = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_100_400_test_fewshot_correct.txt"
$row12[0,45] = 1323.01
$row12[0,46] = 0.94375
$row12[0,47] = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_100_400_test_fewshot_executable.txt"
$row12[0,48] = 0
$ws.Range("A12:AW12").Value = $row12

$row13 = New-Object 'object[,]' 1,49
$row13[0,0] = "llama3:8b-instruct-fp16"
$row13[0,1] = "llama3:70b"
$row13[0,2] = 100
$row13[0,3] = 400
$row13[0,4] = 2657.96
$row13[0,5] = 1317.85
$row13[0,6] = 0.3
$row13[0,7] = "logs\llama3_8b_instruct_fp16_llama3_70b_100_400_test_match.txt"
$row13[0,8] = 1317.85
$row13[0,9] = 0.925
$row13[0,10] = "logs\llama3_8b_instruct_fp16_llama3_70b_100_400_test_correct.txt"
$row13[0,11] = 1317.85
$row13[0,12] = 0.9375
$row13[0,13] = "logs\llama3_8b_instruct_fp16_llama3_70b_100_400_test_executable.txt"
$row13[0,14] = 0
$row13[0,28] = "text"
$row13[0,29] = 140
$row13[0,30] = 0.1
$row13[0,31] = 150
$row13[0,32] = 0.9
$row13[0,33] = 5
$row13[0,34] = 1
$row13[0,35] = 1.2
$row13[0,36] = 1
$row13[0,37] = 1024
$row13[0,38] = 0
$row13[0,39] = 1340.11
$row13[0,40] = 0.2625
$row13[0,41] = "logs\llama3_8b_instruct_fp16_llama3_70b_100_400_test_fewshot_match.txt"
$row13[0,42] = 1340.11
$row13[0,43] = 0.84375
$row13[0,44] = "logs\llama3_8b_instruct_fp16_llama3_70b_100_400_test_fewshot_correct.txt"
$row13[0,45] = 1340.11
$row13[0,46] = 0.93125
$row13[0,47] = "logs\llama3_8b_instruct_fp16_llama3_70b_100_400_test_fewshot_executable.txt"
$row13[0,48] = 0
$ws.Range("A13:AW13").Value = $row13

$row14 = New-Object 'object[,]' 1,49
$row14[0,0] = "codegemma:7b-code-fp16"
$row14[0,1] = "llama3:70b"
$row14[0,2] = 100
$row14[0,3] = 400
$row14[0,4] = 3991.24
$row14[0,5] = 1958.43
$row14[0,6] = 0.23125
$row14[0,7] = "logs\codegemma_7b_code_fp16_llama3_70b_100_400_test_match.txt"
$row14[0,8] = 1958.43
$row14[0,9] = 0.68125
$row14[0,10] = "logs\codegemma_7b_code_fp16_llama3_70b_100_400_test_correct.txt"
$row14[0,11] = 1958.43
$row14[0,12] = 0.425
$row14[0,13] = "logs\codegemma_7b_code_fp16_llama3_70b_100_400_test_executable.txt"
$row14[0,14] = 0
$row14[0,28] = "text"
$row14[0,29] = 140
$row14[0,30] = 0.1
$row14[0,31] = 150
$row14[0,32] = 0.9
$row14[0,33] = 5
$row14[0,34] = 1
$row14[0,35] = 1.2
$row14[0,36] = 1
$row14[0,37] = 1024
$row14[0,38] = 0
$row14[0,39] = 2032.81
$row14[0,40] = 0.33125
$row14[0,41] = "logs\codegemma_7b_code_fp16_llama3_70b_100_400_test_fewshot_match.txt"
$row14[0,42] = 2032.81
$row14[0,43] = 0.88125
$row14[0,44] = "logs\codegemma_7b_code_fp16_llama3_70b_100_400_test_fewshot_correct.txt"
$row14[0,45] = 2032.81
$row14[0,46] = 0.25
$row14[0,47] = "logs\codegemma_7b_code_fp16_llama3_70b_100_400_test_fewshot_executable.txt"
$row14[0,48] = 0
$ws.Range("A14:AW14").Value = $row14
